# Weekly refresh of the Rabanito (radish) price series:
# a new weekly observation is inserted right before the existing row 435
# (old row ref D435/J435), pushing every subsequent record (old rows
# 435-520) down by one row, so the sheet grows from 520 to 521 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 435; Excel shifts rows 435..520 down to
# 436..521 and extends the sheet dimension to A1:R521 automatically.
$ws.Rows("435").Insert()

# Populate the newly inserted row 435 with the new weekly observation.
$ws.Cells.Item(435, 1).Value = 9
$ws.Cells.Item(435, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(435, 3).Value = "Metropolitana"
$ws.Cells.Item(435, 4).Value = 45258
$ws.Cells.Item(435, 5).Value = 13
$ws.Cells.Item(435, 6).Value = 300000001
$ws.Cells.Item(435, 7).Value = "Rabanito"
$ws.Cells.Item(435, 8).Value = "Sin especificar"
$ws.Cells.Item(435, 9).Value = "Primera"
$ws.Cells.Item(435, 10).Value = 8000
$ws.Cells.Item(435, 11).Value = 3000
$ws.Cells.Item(435, 12).Value = 3000
$ws.Cells.Item(435, 13).Value = 3000
$ws.Cells.Item(435, 14).Value = "`$/cien unidades (volumen en unidades)"
$ws.Cells.Item(435, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(435, 16).Value = 30
$ws.Cells.Item(435, 17).Value = 100
$ws.Cells.Item(435, 18).Value = "Hortaliza"
